$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9463
$ws1.Range("F5").Value = 522
$ws1.Range("F6").Value = 463

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9463
$ws4.Range("F5").Value = 522
$ws4.Range("F7").Value = 463
